$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- "core_principles" row (row 8): mark started, log word count ---
$ws.Range("C8").Value = 30
$ws.Range("E8").Value = "y"

# --- "values_goals" row (row 9) dropped from the log ---
$ws.Range("B9").ClearContents()

# --- "open" row (row 11): mark started/finished, log word count ---
$ws.Range("C11").Value = 124
$ws.Range("E11").Value = "y"
$ws.Range("F11").Value = "y"

# --- new sections written up: insert 3 fresh rows right after "open" ---
$ws.Rows("12:14").Insert()

# Populate the new rows (write "open science " before "software" so the
# shared-string table keeps the same insertion order as the authored edit)
$ws.Range("B13").Value = "open science "
$ws.Range("B12").Value = "software"
$ws.Range("B14").Value = "oer"

$ws.Range("C12").Value = 854
$ws.Range("E12").Value = "y"
$ws.Range("F12").Value = "y"

$ws.Range("C13").Value = 3552
$ws.Range("E13").Value = "y"
$ws.Range("F13").Value = "y"

$ws.Range("C14").Value = 1449
$ws.Range("E14").Value = "y"
$ws.Range("F14").Value = "y"

# Restore the "section divider" look (explicit black font) on the blank
# row that now follows the new "oer" entry, matching row 14/15 styling.
$ws.Range("A14").Font.Color = 0
$ws.Range("A15").Font.Color = 0

# --- reflect the latest cursor position left in the workbook ---
[void]$ws.Range("C15").Select()

Write-Host "done"
